$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.589.41'
$ws.Range("E2").Value = '  -1.03%  '

$ws.Range("D3").Value = '2.223.47'
$ws.Range("E3").Value = '  -1.89%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.36%  '

$ws.Range("E7").Value = '  -3.10%  '

$ws.Range("E8").Value = '  +0.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.97'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0820'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.89%  '

$ws.Range("E13").Value = '  -3.06%  '

$ws.Range("D14").Value = '2.562.22'
$ws.Range("E14").Value = '  -1.93%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.260.64'
$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.837'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.26%  '

$ws.Range("D18").Value = '43.485.81'
$ws.Range("E18").Value = '  -0.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -10.09%  '

$ws.Range("D20").Value = '0.0₃0962'
$ws.Range("E20").Value = '  -2.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.48%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.86%  '

$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.40%  '

$ws.Range("E25").Value = '  -7.41%  '

$ws.Range("E26").Value = '  +0.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.10'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.83%  '

$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '160.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.34%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.87'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0827'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.11%  '

$ws.Range("E34").Value = '  -1.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.18'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.89'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.108'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.23%  '

$ws.Range("E38").Value = '  -3.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.56'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -12.69%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0306'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.79%  '

$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").Value = '1.706.87'
$ws.Range("E44").Value = '  -4.05%  '

$ws.Range("B45").Value = 'BitcoinSV'
$ws.Range("C45").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '81.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.64%  '

$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.194'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.53%  '

$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.59%  '

$ws.Range("B48").Value = 'ordi'
$ws.Range("C48").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '72.87'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.68'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.64'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.42'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.25%  '
